{"js": "// Replace the \"Summary 2\" / \"N 2\" column headers in the results table\n// with \"Summary 1\" / \"N 0\" respectively (pt_base v1.1.0 doc update).\n\n// --- \"Summary 2\" -> \"Summary 1\" -----------------------------------------\nlet summaryResults = context.document.body.search(\"Summary 2\", { matchCase: true });\nsummaryResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < summaryResults.items.length; i++) {\n  summaryResults.items[i].insertText(\"Summary 1\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- \"N 2\" -> \"N 0\" -------------------------------------------------------\nlet nResults = context.document.body.search(\"N 2\", { matchCase: true });\nnResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < nResults.items.length; i++) {\n  nResults.items[i].insertText(\"N 0\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Replace the \"Summary 2\" / \"N 2\" column headers in the results table\n# with \"Summary 1\" / \"N 0\" respectively (pt_base v1.1.0 doc update).\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n# Header row (row 1): ... | Summary 2 (col 5) | N 2 (col 6) | ...\n$table.Cell(1, 5).Range.Text = \"Summary 1\"\n$table.Cell(1, 6).Range.Text = \"N 0\"\n"}
